$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.548.63'
$ws.Range('E2').Value = '  +2.56%  '
$ws.Range('D3').Value = '1.641.02'
$ws.Range('E3').Value = '  +4.24%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9997'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '308.36'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.01%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9988'
$ws.Range('D6').Style = "Normal"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3782'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.81%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '53.11'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +6.41%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3689'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +3.85%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.284'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +5.96%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08212'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +3.27%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.9986'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.38%  '
$ws.Range('E13').Value = '  +7.02%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.679'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +4.64%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.00001291'
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '7.490'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.89%  '
$ws.Range('D17').Value = '1.641.53'
$ws.Range('E17').Value = '  +4.22%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '95.04'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +3.45%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06958'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +3.31%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.46'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +4.76%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.602'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +4.21%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9992'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('D23').Value = '23.543.76'
$ws.Range('E23').Value = '  +2.60%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.01'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +3.44%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.157'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +12.38%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.414'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.73%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '21.47'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +4.62%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '151.78'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.51%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.348'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +3.53%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '136.54'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +4.19%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.428'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +4.13%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.861'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +5.01%  '
$ws.Range('D33').Value = '1.817.68'
$ws.Range('E33').Value = '  +3.87%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.9786'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +5.35%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.02819'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +7.33%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '10.46'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +5.48%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.07499'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.30%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.242'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +4.87%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.2548'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +3.85%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.08870'
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.402'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +4.48%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.7200'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +5.19%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '12.68'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +7.48%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.15'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +9.55%  '
$ws.Range('E45').Value = '  +5.62%  '
$ws.Range('E46').Value = '  +5.87%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.047'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.19%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.9979'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.31%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.08073'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.00%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '131.50'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.63%  '
$ws.Range('E51').Value = '  +3.40%  '
